# Add the new "Notes" entries (screw type / thread callouts) for the
# hardware rows, as part of improving the BOM with O-rings / module PCB info.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = "8-32x5/8 flat"
$ws.Range("F5").Value = "4-40x1/8"
$ws.Range("F6").Value = "4-40x5/16"
$ws.Range("F7").Value = "4-40x3/8 torx"

# Re-size the McMaster Carr (D) and new Notes (F) columns to fit their
# (now much shorter) contents.
$ws.Columns.Item(4).ColumnWidth = 12.8
$ws.Columns.Item(6).ColumnWidth = 11

# Leave the selection on the last-edited cell.
$ws.Range("F8").Select() | Out-Null
